# Scenario_HeatingTechnology_Availability.xlsx
#
# Insert a new leading column "id_scenario" (constant value 1 for every
# data row) in front of "id_region" on the single worksheet, shifting the
# existing id_region / id_heating_technology / id_heating_system_action /
# unit / year columns one place to the right. The table (Table1) is
# rebuilt over the new A1:AT76 range so its column list / header names
# stay in sync with the shifted header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- locate the existing table and remember its used range -----------
$tbl = $ws.ListObjects.Item(1)
$lastRow = $tbl.Range.Rows.Count               # 76 (incl. header)
$lastCol = $tbl.Range.Columns.Count            # 45 (A..AS)

# Convert the table back to a plain range first: in this environment
# ListObject.Delete() removes the underlying cell data as well, whereas
# Unlist() keeps all the cell values/headers intact and just drops the
# table wrapper, which is what we need before reshaping the grid.
$tbl.Unlist()

# --- insert the new column and populate it ----------------------------
$ws.Range("A1").EntireColumn.Insert()
$ws.Range("A1").Value2 = "id_scenario"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 1
}

# --- rebuild the table over the new range ------------------------------
$newLastCol = $lastCol + 1                      # 46 (A..AT)
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $newLastCol))
$newTbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$newTbl.Name = "Table1"
$newTbl.TableStyle = "TableStyleMedium6"

# --- cosmetic: restore a sensible selection ----------------------------
$ws.Range("F24").Select()
